$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Training Dashboard")

# --- Header styling: the dashboard title and the column-header row both
# end up rendered as bold white text (the title drops its old explicit
# 14pt size so it shares the same "bold, white" look as the blue header
# band). Touching both ranges in the same pass lets them collapse onto a
# single shared font definition.
$ws.Range("A1").Font.Bold = $true
$ws.Range("A1").Font.Size = 11
$ws.Range("A1").Font.Color = 16777215

$ws.Range("A2:K2").Font.Bold = $true
$ws.Range("A2:K2").Font.Color = 16777215

# --- Refresh the "PERIOD TO EXPIRE" / "LAST UPDATE" columns for the
# three training rows (new run pulled on 16-Sep-2025 instead of
# 08-Sep-2025, shifting the day counts by -8). The leading apostrophe
# keeps the date column as literal text instead of Excel auto-converting
# it to a date serial number.
$ws.Range("H3").Value = 583
$ws.Range("I3").Value = "'16-Sep-2025"

$ws.Range("H4").Value = 583
$ws.Range("I4").Value = "'16-Sep-2025"

$ws.Range("H5").Value = -19618
$ws.Range("I5").Value = "'16-Sep-2025"
